$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 31-32, shifting existing rows 31-41 down to 33-43.
[void]$ws.Range("A31:R32").Insert()

# --- New row 31 (Primera, $/paquete 6 unidades, Region de Nuble) ---
[void]($ws.Cells.Item(31,1).Value() = 11)
[void]($ws.Cells.Item(31,2).Value() = 'Vega Monumental Concepción')
[void]($ws.Cells.Item(31,3).Value() = 'Bíobío')
[void]($ws.Cells.Item(31,4).Value() = 44694)
[void]($ws.Cells.Item(31,5).Value() = 8)
[void]($ws.Cells.Item(31,6).Value() = 100112037)
[void]($ws.Cells.Item(31,7).Value() = 'Cebollín')
[void]($ws.Cells.Item(31,8).Value() = 'Sin especificar')
[void]($ws.Cells.Item(31,9).Value() = 'Primera')
[void]($ws.Cells.Item(31,10).Value() = 200)
[void]($ws.Cells.Item(31,11).Value() = 600)
[void]($ws.Cells.Item(31,12).Value() = 700)
[void]($ws.Cells.Item(31,13).Value() = 650)
[void]($ws.Cells.Item(31,14).Value() = '$/paquete 6 unidades')
[void]($ws.Cells.Item(31,15).Value() = 'Región de Ñuble')
[void]($ws.Cells.Item(31,16).Value() = 108)
[void]($ws.Cells.Item(31,17).Value() = 6)
[void]($ws.Cells.Item(31,18).Value() = 'Hortaliza')

# --- New row 32 (Segunda, $/paquete 6 unidades, Region de Nuble) ---
[void]($ws.Cells.Item(32,1).Value() = 11)
[void]($ws.Cells.Item(32,2).Value() = 'Vega Monumental Concepción')
[void]($ws.Cells.Item(32,3).Value() = 'Bíobío')
[void]($ws.Cells.Item(32,4).Value() = 44694)
[void]($ws.Cells.Item(32,5).Value() = 8)
[void]($ws.Cells.Item(32,6).Value() = 100112037)
[void]($ws.Cells.Item(32,7).Value() = 'Cebollín')
[void]($ws.Cells.Item(32,8).Value() = 'Sin especificar')
[void]($ws.Cells.Item(32,9).Value() = 'Segunda')
[void]($ws.Cells.Item(32,10).Value() = 100)
[void]($ws.Cells.Item(32,11).Value() = 500)
[void]($ws.Cells.Item(32,12).Value() = 500)
[void]($ws.Cells.Item(32,13).Value() = 500)
[void]($ws.Cells.Item(32,14).Value() = '$/paquete 6 unidades')
[void]($ws.Cells.Item(32,15).Value() = 'Región de Ñuble')
[void]($ws.Cells.Item(32,16).Value() = 83)
[void]($ws.Cells.Item(32,17).Value() = 6)
[void]($ws.Cells.Item(32,18).Value() = 'Hortaliza')

Write-Output $ws.UsedRange.Address()
